{"js": "const pairs = [\n  [\"53-33=\", \"4+34=\"],\n  [\"58-38=\", \"53+41=\"],\n  [\"93-11=\", \"72+10=\"],\n  [\"64-48=\", \"79-11=\"],\n  [\"54-18=\", \"46-31=\"],\n  [\"76-59=\", \"14+2=\"],\n  [\"70-3=\", \"37-34=\"],\n  [\"83-70=\", \"39+17=\"],\n  [\"3+43=\", \"43+44=\"],\n  [\"56-16=\", \"93-34=\"],\n  [\"77-19=\", \"13-7=\"],\n  [\"86-31=\", \"10+41=\"],\n  [\"96-80=\", \"18-3=\"],\n  [\"88-78=\", \"67-13=\"],\n  [\"13+0=\", \"83-43=\"],\n  [\"92-25=\", \"17+73=\"],\n  [\"51+35=\", \"59-59=\"],\n  [\"14+59=\", \"48+18=\"],\n  [\"34-32=\", \"48-2=\"],\n  [\"71-58=\", \"99-72=\"],\n  [\"0+45=\", \"33+56=\"],\n  [\"74-11=\", \"11+74=\"],\n  [\"12+48=\", \"60-51=\"],\n  [\"84-78=\", \"22+68=\"],\n  [\"68-12=\", \"59-12=\"],\n  [\"76+6=\", \"68-19=\"],\n  [\"79-39=\", \"1+96=\"],\n  [\"0+21=\", \"87-32=\"],\n  [\"64-29=\", \"88-79=\"],\n  [\"23+33=\", \"8+9=\"],\n  [\"98-24=\", \"68-50=\"],\n  [\"43-9=\", \"86-26=\"],\n  [\"65+30=\", \"77-49=\"],\n  [\"51-20=\", \"30+66=\"],\n  [\"15+44=\", \"60+8=\"],\n  [\"74-63=\", \"93-77=\"],\n  [\"61-52=\", \"37+54=\"],\n  [\"10+25=\", \"66-31=\"],\n  [\"9+56=\", \"21+51=\"],\n  [\"30+37=\", \"55-29=\"],\n  [\"43+19=\", \"95-13=\"],\n  [\"22+42=\", \"68-45=\"],\n  [\"52+45=\", \"94-69=\"],\n  [\"87-35=\", \"14+0=\"],\n  [\"80+16=\", \"13+15=\"],\n  [\"52+44=\", \"84-62=\"],\n  [\"24+57=\", \"29+57=\"],\n  [\"27+51=\", \"19+79=\"],\n  [\"46+28=\", \"20+43=\"],\n  [\"55-26=\", \"22+52=\"],\n  [\"46+33=\", \"59-0=\"],\n  [\"24+49=\", \"34+0=\"],\n  [\"74-13=\", \"41-8=\"],\n  [\"12+29=\", \"57-45=\"],\n  [\"68-24=\", \"94-72=\"],\n  [\"71-30=\", \"27+26=\"],\n  [\"61-13=\", \"5+50=\"],\n  [\"70-58=\", \"67+7=\"],\n  [\"11+1=\", \"48-31=\"],\n  [\"31-14=\", \"61+36=\"],\n  [\"12+50=\", \"23-19=\"],\n  [\"40+9=\", \"3+19=\"],\n  [\"61+26=\", \"78+5=\"],\n  [\"82+3=\", \"61-45=\"],\n  [\"41-21=\", \"1+13=\"],\n  [\"65+33=\", \"58-47=\"],\n  [\"71+24=\", \"54-9=\"],\n  [\"52+42=\", \"44-28=\"],\n  [\"28-5=\", \"72-62=\"],\n  [\"37+14=\", \"77-32=\"],\n  [\"40+20=\", \"90-35=\"],\n  [\"57-9=\", \"73-67=\"],\n  [\"78+7=\", \"67-16=\"],\n  [\"58+11=\", \"19+60=\"],\n  [\"37+44=\", \"1+50=\"],\n  [\"9+47=\", \"97-86=\"],\n  [\"55-51=\", \"72-36=\"],\n  [\"48+13=\", \"44+8=\"],\n  [\"39+27=\", \"20+22=\"],\n  [\"97+2=\", \"40-3=\"],\n  [\"3+38=\", \"86-70=\"],\n  [\"20-9=\", \"56-9=\"],\n  [\"15+25=\", \"25+67=\"],\n  [\"97-16=\", \"43+36=\"],\n  [\"19+69=\", \"12+9=\"],\n  [\"80-61=\", \"28-18=\"],\n  [\"2-2=\", \"4+22=\"],\n  [\"78-43=\", \"18+39=\"],\n  [\"5+49=\", \"66-17=\"],\n  [\"64+26=\", \"0+26=\"],\n  [\"57-0=\", \"15+6=\"],\n  [\"3+37=\", \"47-19=\"],\n  [\"28-19=\", \"54-19=\"],\n  [\"2+69=\", \"95-75=\"],\n  [\"5+61=\", \"9+87=\"],\n  [\"40+29=\", \"39-23=\"],\n  [\"78-35=\", \"6+93=\"],\n  [\"60-46=\", \"62+1=\"],\n  [\"82-6=\", \"91-34=\"],\n  [\"70-32=\", \"70-22=\"],\n];\n\nconst body = context.document.body;\n\n// Phase 1: issue a search for every old equation text.\nconst searches = pairs.map(([oldText, newText]) => {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  return { results, oldText, newText };\n});\nawait context.sync();\n\n// Phase 2: replace every match found above with its new equation text.\nfor (const { results, newText } of searches) {\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"53-33=\", \"4+34=\"),\n    @(\"58-38=\", \"53+41=\"),\n    @(\"93-11=\", \"72+10=\"),\n    @(\"64-48=\", \"79-11=\"),\n    @(\"54-18=\", \"46-31=\"),\n    @(\"76-59=\", \"14+2=\"),\n    @(\"70-3=\", \"37-34=\"),\n    @(\"83-70=\", \"39+17=\"),\n    @(\"3+43=\", \"43+44=\"),\n    @(\"56-16=\", \"93-34=\"),\n    @(\"77-19=\", \"13-7=\"),\n    @(\"86-31=\", \"10+41=\"),\n    @(\"96-80=\", \"18-3=\"),\n    @(\"88-78=\", \"67-13=\"),\n    @(\"13+0=\", \"83-43=\"),\n    @(\"92-25=\", \"17+73=\"),\n    @(\"51+35=\", \"59-59=\"),\n    @(\"14+59=\", \"48+18=\"),\n    @(\"34-32=\", \"48-2=\"),\n    @(\"71-58=\", \"99-72=\"),\n    @(\"0+45=\", \"33+56=\"),\n    @(\"74-11=\", \"11+74=\"),\n    @(\"12+48=\", \"60-51=\"),\n    @(\"84-78=\", \"22+68=\"),\n    @(\"68-12=\", \"59-12=\"),\n    @(\"76+6=\", \"68-19=\"),\n    @(\"79-39=\", \"1+96=\"),\n    @(\"0+21=\", \"87-32=\"),\n    @(\"64-29=\", \"88-79=\"),\n    @(\"23+33=\", \"8+9=\"),\n    @(\"98-24=\", \"68-50=\"),\n    @(\"43-9=\", \"86-26=\"),\n    @(\"65+30=\", \"77-49=\"),\n    @(\"51-20=\", \"30+66=\"),\n    @(\"15+44=\", \"60+8=\"),\n    @(\"74-63=\", \"93-77=\"),\n    @(\"61-52=\", \"37+54=\"),\n    @(\"10+25=\", \"66-31=\"),\n    @(\"9+56=\", \"21+51=\"),\n    @(\"30+37=\", \"55-29=\"),\n    @(\"43+19=\", \"95-13=\"),\n    @(\"22+42=\", \"68-45=\"),\n    @(\"52+45=\", \"94-69=\"),\n    @(\"87-35=\", \"14+0=\"),\n    @(\"80+16=\", \"13+15=\"),\n    @(\"52+44=\", \"84-62=\"),\n    @(\"24+57=\", \"29+57=\"),\n    @(\"27+51=\", \"19+79=\"),\n    @(\"46+28=\", \"20+43=\"),\n    @(\"55-26=\", \"22+52=\"),\n    @(\"46+33=\", \"59-0=\"),\n    @(\"24+49=\", \"34+0=\"),\n    @(\"74-13=\", \"41-8=\"),\n    @(\"12+29=\", \"57-45=\"),\n    @(\"68-24=\", \"94-72=\"),\n    @(\"71-30=\", \"27+26=\"),\n    @(\"61-13=\", \"5+50=\"),\n    @(\"70-58=\", \"67+7=\"),\n    @(\"11+1=\", \"48-31=\"),\n    @(\"31-14=\", \"61+36=\"),\n    @(\"12+50=\", \"23-19=\"),\n    @(\"40+9=\", \"3+19=\"),\n    @(\"61+26=\", \"78+5=\"),\n    @(\"82+3=\", \"61-45=\"),\n    @(\"41-21=\", \"1+13=\"),\n    @(\"65+33=\", \"58-47=\"),\n    @(\"71+24=\", \"54-9=\"),\n    @(\"52+42=\", \"44-28=\"),\n    @(\"28-5=\", \"72-62=\"),\n    @(\"37+14=\", \"77-32=\"),\n    @(\"40+20=\", \"90-35=\"),\n    @(\"57-9=\", \"73-67=\"),\n    @(\"78+7=\", \"67-16=\"),\n    @(\"58+11=\", \"19+60=\"),\n    @(\"37+44=\", \"1+50=\"),\n    @(\"9+47=\", \"97-86=\"),\n    @(\"55-51=\", \"72-36=\"),\n    @(\"48+13=\", \"44+8=\"),\n    @(\"39+27=\", \"20+22=\"),\n    @(\"97+2=\", \"40-3=\"),\n    @(\"3+38=\", \"86-70=\"),\n    @(\"20-9=\", \"56-9=\"),\n    @(\"15+25=\", \"25+67=\"),\n    @(\"97-16=\", \"43+36=\"),\n    @(\"19+69=\", \"12+9=\"),\n    @(\"80-61=\", \"28-18=\"),\n    @(\"2-2=\", \"4+22=\"),\n    @(\"78-43=\", \"18+39=\"),\n    @(\"5+49=\", \"66-17=\"),\n    @(\"64+26=\", \"0+26=\"),\n    @(\"57-0=\", \"15+6=\"),\n    @(\"3+37=\", \"47-19=\"),\n    @(\"28-19=\", \"54-19=\"),\n    @(\"2+69=\", \"95-75=\"),\n    @(\"5+61=\", \"9+87=\"),\n    @(\"40+29=\", \"39-23=\"),\n    @(\"78-35=\", \"6+93=\"),\n    @(\"60-46=\", \"62+1=\"),\n    @(\"82-6=\", \"91-34=\"),\n    @(\"70-32=\", \"70-22=\"),\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n\nWrite-Output \"Done\""}
